$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append the new weekly observation as row 95 ---
$dataSheet = $wb.Worksheets.Item("Data")

# Clone the formatting of the previous data row (A94:B94) down into row 95
# so the new date cell keeps the same date number-format/border/font as the
# rest of column A, then overwrite with the new observation's values.
$dataSheet.Range("A94:B94").Copy()
$dataSheet.Range("A95:B95").PasteSpecial(-4122)  # xlPasteFormats
$dataSheet.Range("A95").Value = 45126
$dataSheet.Range("B95").Value = 8274.552

# --- Sheet "SeriesInfo": refresh metadata pulled from the FRED API ---
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

# Force text format first so these date-looking strings are not
# auto-converted into date serial numbers.
$infoSheet.Range("B3").NumberFormat = "@"
$infoSheet.Range("B3").Value = "2023-07-24"

$infoSheet.Range("B4").NumberFormat = "@"
$infoSheet.Range("B4").Value = "2023-07-24"

$infoSheet.Range("B7").NumberFormat = "@"
$infoSheet.Range("B7").Value = "2023-07-19"

$infoSheet.Range("B14").NumberFormat = "@"
$infoSheet.Range("B14").Value = "2023-07-20 15:33:32-05"
